$d = $word.ActiveDocument

# Locate the paragraph containing ".func or .function [String name]"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ".func or .function*") {
        $target = $p
        break
    }
}

# Insert the new ".params or .args or .argcount <Number>" paragraph right
# before the ".func" paragraph.
$target.Range.InsertBefore(".params or .args or .argcount <Number>`r")

# Re-locate the ".func" paragraph (index shifted by the inserted paragraph)
# and give it the lastRenderedPageBreak marker that used to sit on the
# ".end" paragraph.
$funcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ".func or .function*") {
        $funcPara = $p
        break
    }
}
$funcXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>.func or .function [String name]</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> &#8211; begins a new function in the current one</w:t></w:r></w:p>'
$null = $funcPara.Range.InsertXML($funcXml)

# Remove the lastRenderedPageBreak from the ".end" paragraph, since it has
# moved to the ".func" paragraph above.
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ".end*") {
        $endPara = $p
        break
    }
}
$endXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009024DC" w:rsidRDefault="009024DC" w:rsidP="00456A0C"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.end &#8211; ends a function</w:t></w:r></w:p>'
$null = $endPara.Range.InsertXML($endXml)

Write-Output "Edit applied"
